# ===================================================================
# Edit: add a new '2022-Q4' sheet (inserted before '2022-Q3') with the
# quarter's fund-holding detail, and add the corresponding summary row
# to the '总计' (totals) sheet.
# ===================================================================

$wb = $excel.ActiveWorkbook

# --- 1. Update '总计' (totals) sheet: insert a new row 2 for 2022-Q4 ---
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()

# Copy formatting from the (now-shifted) existing data rows so the new
# row matches the sheet's existing look (bold/centered/bordered index
# column, plain data columns).
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$wsTotal.Range("B3:D3").Copy()
$wsTotal.Range("B2:D2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 12
$wsTotal.Range("D2").Value = 0.48

# Re-sequence the index column for the rows that shifted down.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2

# --- 2. Insert the new '2022-Q4' worksheet, positioned right after
#        '总计' and before '2022-Q3'. The most reliable way to get an
#        exact structural/formatting clone (sheetPr, header style,
#        index-column style, …) in this host is to duplicate the
#        existing '2022-Q3' sheet and then overwrite its values ---
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($wsQ3)
$wsNew = $wb.Worksheets.Item("2022-Q3 (2)")
$wsNew.Name = "2022-Q4"

# '2022-Q3' only has 2 data rows (3 rows total); 2022-Q4 needs 12 data
# rows (13 rows total). Extend the data area downward, cloning row 3's
# formatting (bold/bordered index column, plain data columns) so every
# new row matches the sheet's existing look.
$wsNew.Range("A3:H3").Copy()
$wsNew.Range("A4:H13").PasteSpecial(-4122)

# Header row (already holds the right labels/style from the sheet
# copy -- re-asserted here only for robustness/clarity).
$wsNew.Range("B1").Value = "基金代码"
$wsNew.Range("C1").Value = "基金名称"
$wsNew.Range("D1").Value = "基金规模"
$wsNew.Range("E1").Value = "股票总仓位"
$wsNew.Range("F1").Value = "仓位占比"
$wsNew.Range("G1").Value = "持有市值(亿元)"
$wsNew.Range("H1").Value = "仓位排名"

# Data rows. A = running index, H = rank -> genuine numbers.
# B (fund code) and D/E/F/G (scale/position/ratio/market-value) are
# stored as literal text in the source data (e.g. fund codes keep
# leading zeros, market values keep trailing zeros like "0.0780"),
# so they are written with a leading apostrophe -- the standard
# Excel convention for entering a numeric-looking value as text --
# which keeps the exact string instead of it being parsed as a number.
# --- row 2 ---
$wsNew.Range("A2").Value = 0
$wsNew.Range("B2").Value = "'013676"
$wsNew.Range("C2").Value = "兴银兴慧一年持有混合A"
$wsNew.Range("D2").Value = "'8.13"
$wsNew.Range("E2").Value = "'23.86"
$wsNew.Range("F2").Value = "'0.96"
$wsNew.Range("G2").Value = "'0.0780"
$wsNew.Range("H2").Value = 7
# --- row 3 ---
$wsNew.Range("A3").Value = 1
$wsNew.Range("B3").Value = "'004818"
$wsNew.Range("C3").Value = "国寿安保目标策略灵活配置混合A"
$wsNew.Range("D3").Value = "'2.76"
$wsNew.Range("E3").Value = "'59.92"
$wsNew.Range("F3").Value = "'2.51"
$wsNew.Range("G3").Value = "'0.0693"
$wsNew.Range("H3").Value = 10
# --- row 4 ---
$wsNew.Range("A4").Value = 2
$wsNew.Range("B4").Value = "'016588"
$wsNew.Range("C4").Value = "富国融甄混合A"
$wsNew.Range("D4").Value = "'4.57"
$wsNew.Range("E4").Value = "'29.63"
$wsNew.Range("F4").Value = "'1.50"
$wsNew.Range("G4").Value = "'0.0686"
$wsNew.Range("H4").Value = 8
# --- row 5 ---
$wsNew.Range("A5").Value = 3
$wsNew.Range("B5").Value = "'014839"
$wsNew.Range("C5").Value = "兴银碳中和主题混合C"
$wsNew.Range("D5").Value = "'0.64"
$wsNew.Range("E5").Value = "'92.17"
$wsNew.Range("F5").Value = "'8.21"
$wsNew.Range("G5").Value = "'0.0525"
$wsNew.Range("H5").Value = 1
# --- row 6 ---
$wsNew.Range("A6").Value = 4
$wsNew.Range("B6").Value = "'013677"
$wsNew.Range("C6").Value = "兴银兴慧一年持有混合C"
$wsNew.Range("D6").Value = "'4.59"
$wsNew.Range("E6").Value = "'23.86"
$wsNew.Range("F6").Value = "'0.96"
$wsNew.Range("G6").Value = "'0.0441"
$wsNew.Range("H6").Value = 7
# --- row 7 ---
$wsNew.Range("A7").Value = 5
$wsNew.Range("B7").Value = "'014838"
$wsNew.Range("C7").Value = "兴银碳中和主题混合A"
$wsNew.Range("D7").Value = "'0.53"
$wsNew.Range("E7").Value = "'92.17"
$wsNew.Range("F7").Value = "'8.21"
$wsNew.Range("G7").Value = "'0.0435"
$wsNew.Range("H7").Value = 1
# --- row 8 ---
$wsNew.Range("A8").Value = 6
$wsNew.Range("B8").Value = "'004819"
$wsNew.Range("C8").Value = "国寿安保目标策略灵活配置混合C"
$wsNew.Range("D8").Value = "'1.30"
$wsNew.Range("E8").Value = "'59.92"
$wsNew.Range("F8").Value = "'2.51"
$wsNew.Range("G8").Value = "'0.0326"
$wsNew.Range("H8").Value = 10
# --- row 9 ---
$wsNew.Range("A9").Value = 7
$wsNew.Range("B9").Value = "'016589"
$wsNew.Range("C9").Value = "富国融甄混合C"
$wsNew.Range("D9").Value = "'1.96"
$wsNew.Range("E9").Value = "'29.63"
$wsNew.Range("F9").Value = "'1.50"
$wsNew.Range("G9").Value = "'0.0294"
$wsNew.Range("H9").Value = 8
# --- row 10 ---
$wsNew.Range("A10").Value = 8
$wsNew.Range("B10").Value = "'010540"
$wsNew.Range("C10").Value = "浙商智多金稳健一年持有期混合C"
$wsNew.Range("D10").Value = "'1.37"
$wsNew.Range("E10").Value = "'25.01"
$wsNew.Range("F10").Value = "'1.33"
$wsNew.Range("G10").Value = "'0.0182"
$wsNew.Range("H10").Value = 7
# --- row 11 ---
$wsNew.Range("A11").Value = 9
$wsNew.Range("B11").Value = "'010539"
$wsNew.Range("C11").Value = "浙商智多金稳健一年持有期混合A"
$wsNew.Range("D11").Value = "'1.27"
$wsNew.Range("E11").Value = "'25.01"
$wsNew.Range("F11").Value = "'1.33"
$wsNew.Range("G11").Value = "'0.0169"
$wsNew.Range("H11").Value = 7
# --- row 12 ---
$wsNew.Range("A12").Value = 10
$wsNew.Range("B12").Value = "'009569"
$wsNew.Range("C12").Value = "浙商智多宝稳健一年持有期混合C"
$wsNew.Range("D12").Value = "'1.02"
$wsNew.Range("E12").Value = "'26.91"
$wsNew.Range("F12").Value = "'1.20"
$wsNew.Range("G12").Value = "'0.0122"
$wsNew.Range("H12").Value = 4
# --- row 13 ---
$wsNew.Range("A13").Value = 11
$wsNew.Range("B13").Value = "'009568"
$wsNew.Range("C13").Value = "浙商智多宝稳健一年持有期混合A"
$wsNew.Range("D13").Value = "'0.98"
$wsNew.Range("E13").Value = "'26.91"
$wsNew.Range("F13").Value = "'1.20"
$wsNew.Range("G13").Value = "'0.0118"
$wsNew.Range("H13").Value = 4

Write-Host "edit complete"
